# Update Name of Algo
# Apply updated imputed values to specific cells on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 6.633300000000004
$ws.Range("B21").Value = 9.1882
$ws.Range("B23").Value = 8.900799999999998
$ws.Range("D24").Value = -7.613299999999997
$ws.Range("B25").Value = 5.786299999999998
$ws.Range("D28").Value = -8.345799999999999
$ws.Range("D36").Value = -7.1096
$ws.Range("D45").Value = -7.227600000000001
$ws.Range("D48").Value = -7.451099999999994
$ws.Range("D49").Value = -7.974999999999998
$ws.Range("D52").Value = -8.032000000000007
$ws.Range("B53").Value = 5.160900000000002
$ws.Range("D53").Value = -8.0265
$ws.Range("D54").Value = -8.130600000000006
$ws.Range("B57").Value = 5.006899999999995
$ws.Range("B59").Value = 4.954999999999999
$ws.Range("B69").Value = 5.233499999999996
$ws.Range("D70").Value = -7.1727
$ws.Range("B79").Value = 9.482800000000005
$ws.Range("B83").Value = 4.898999999999998
$ws.Range("D86").Value = -8.286699999999996
$ws.Range("D87").Value = -8.590599999999995
$ws.Range("B93").Value = 5.661599999999998
$ws.Range("D101").Value = -7.832100000000004
